# Fruta / hortaliza, semanal
# The underlying data rows (2-16) get shuffled/reordered as the weekly
# snapshot rolls forward. Column A:T for each destination row comes from a
# (possibly different) source row in the original layout.
#
# Mapping of new-row -> old-row (1-indexed sheet rows):
#   2 <- 4     3 <- 3     4 <- 6     5 <- 2     6 <- 13
#   7 <- 15    8 <- 8     9 <- 14   10 <- 5    11 <- 10
#  12 <- 11   13 <- 12   14 <- 9    15 <- 7    16 <- 16

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastCol = 20   # column T

# Snapshot the current values of every data row (2-16) before overwriting
# anything, since several destinations read from rows that are themselves
# about to be rewritten.
$snapshot = @{}
for ($r = 2; $r -le 16; $r++) {
    $rowVals = @()
    for ($c = 1; $c -le $lastCol; $c++) {
        $rowVals += ,$ws.Cells.Item($r, $c).Value2
    }
    $snapshot[$r] = $rowVals
}

$rowMap = @{
    2  = 4
    3  = 3
    4  = 6
    5  = 2
    6  = 13
    7  = 15
    8  = 8
    9  = 14
    10 = 5
    11 = 10
    12 = 11
    13 = 12
    14 = 9
    15 = 7
    16 = 16
}

foreach ($destRow in $rowMap.Keys) {
    $srcRow = $rowMap[$destRow]
    $vals = $snapshot[$srcRow]
    for ($c = 1; $c -le $lastCol; $c++) {
        $ws.Cells.Item($destRow, $c).Value = $vals[$c - 1]
    }
}
